$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 180, pushing the existing rows 180:225 down to 181:226
$ws.Rows.Item(180).Insert()

# Populate the newly inserted row 180 with the new data record
$ws.Cells.Item(180, 1).Value = 5
$ws.Cells.Item(180, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(180, 3).Value = "Maule"
$ws.Cells.Item(180, 4).Value = 44642
$ws.Cells.Item(180, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(180, 5).Value = 7
$ws.Cells.Item(180, 6).Value = "Fruta"
$ws.Cells.Item(180, 7).Value = 100108
$ws.Cells.Item(180, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(180, 9).Value = 100108005
$ws.Cells.Item(180, 10).Value = "Piña"
$ws.Cells.Item(180, 11).Value = "Caramelo"
$ws.Cells.Item(180, 12).Value = "Segunda"
$ws.Cells.Item(180, 13).Value = 350
$ws.Cells.Item(180, 14).Value = 16000
$ws.Cells.Item(180, 15).Value = 16000
$ws.Cells.Item(180, 16).Value = 16000
$ws.Cells.Item(180, 17).Value = "$/caja 14 unidades"
$ws.Cells.Item(180, 18).Value = "Ecuador"
$ws.Cells.Item(180, 19).Value = 1143
$ws.Cells.Item(180, 20).Value = 14
